# Xor_Genetic.xlsx - network topology grew (2x2 -> 4x2 neurons) and one more
# input/bias row was added; weights were retrained (new random values).
# "params": m (C3) 2 -> 4.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)   # "weights"
$ws2 = $wb.Worksheets.Item(2)   # "params"

# ---------------------------------------------------------------------------
# weights sheet - header rows (1: "Layer:" groups, 2: "Neuron:" indices)
# ---------------------------------------------------------------------------
# Layer 1 now spans 4 neurons (B:E) instead of 2 (B:C); layer 2 still has 2
# neurons but moves from D:E to F:G.
$ws1.Range("B1:E1").Merge()
$ws1.Range("A2").Copy()
$ws1.Range("B1:E1").PasteSpecial(-4122)   # xlPasteFormats - restore the plain header look after Merge's auto border split

$ws1.Range("F1:G1").Merge()
$ws1.Range("A2").Copy()
$ws1.Range("F1:G1").PasteSpecial(-4122)

$ws1.Range("A2").Copy()
$ws1.Range("F2:G2").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws1.Range("B1").Value = 1
$ws1.Range("F1").Value = 2

$ws1.Range("B2").Value = 0
$ws1.Range("C2").Value = 1
$ws1.Range("D2").Value = 2
$ws1.Range("E2").Value = 3
$ws1.Range("F2").Value = 0
$ws1.Range("G2").Value = 1

# ---------------------------------------------------------------------------
# weights sheet - weight matrix, rows 4-6 updated + new layer-2 columns F,G
# (row 3 stays blank, same as before the edit)
# ---------------------------------------------------------------------------

$ws1.Range("B4").Value = 3.995910303355823
$ws1.Range("C4").Value = -9.595948218959755
$ws1.Range("D4").Value = -9.824522234801375
$ws1.Range("E4").Value = 3.012523386384068
$ws1.Range("F4").Value = 16.94989531332135
$ws1.Range("G4").Value = -17.37845937967623

$ws1.Range("B5").Value = 2.0499316506778
$ws1.Range("C5").Value = -4.492484374931386
$ws1.Range("D5").Value = -5.290784327382515
$ws1.Range("E5").Value = 8.555680557054506
$ws1.Range("F5").Value = 6.542476586241063
$ws1.Range("G5").Value = -3.536985016217162

$ws1.Range("B6").Value = -3.721699685396691
$ws1.Range("C6").Value = -5.811602190161497
$ws1.Range("D6").Value = -5.21462714713974
$ws1.Range("E6").Value = 0.0397458734780809
$ws1.Range("F6").Value = 9.501351289512513
$ws1.Range("G6").Value = -17.68165301722782

# ---------------------------------------------------------------------------
# weights sheet - two new input/bias rows (7, 8); only the layer-2 (F,G)
# weights apply to them, B:E are left blank like the rest of the table.
# ---------------------------------------------------------------------------

$ws1.Range("A6").Copy()
$ws1.Range("A7:A8").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws1.Range("A7").Value = 3
$ws1.Range("F7").Value = 2.7888057487887
$ws1.Range("G7").Value = -8.389815711943996

$ws1.Range("A8").Value = 4
$ws1.Range("F8").Value = 7.135461618165352
$ws1.Range("G8").Value = -10.14564680090852

# ---------------------------------------------------------------------------
# params sheet - m (C3) 2 -> 4
# ---------------------------------------------------------------------------

$ws2.Range("C3").Value = 4
